# Update global_variable_generator to include pump variables
#
# - Adds a new "Pump" worksheet (after "Sensor Data") populated with the
#   pump_e_* / pump_h_* global-variable rows.
# - "Shelf" (the previously active tab) loses its tabSelected/topLeftCell
#   view state and gets a new selection (D6:E6); the new "Pump" sheet
#   becomes the selected/active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "Shelf" sheet's lingering view state (it was the
#    previously-active tab) before we move on.
# ---------------------------------------------------------------------
$shelf = $wb.Worksheets.Item("Shelf")
$shelf.Range("D6:E6").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Pump" sheet at the end of the workbook (after
#    "Sensor Data").
# ---------------------------------------------------------------------
$sensorData = $wb.Worksheets.Item($wb.Worksheets.Count)
$pump = $wb.Worksheets.Add([Type]::Missing, $sensorData)
$pump.Name = "Pump"

# Header row
$pump.Range("A1").Value = "base_addr"
$pump.Range("B1").Value = "variable_name"
$pump.Range("C1").Value = "addr_offset"
$pump.Range("D1").Value = "type"
$pump.Range("E1").Value = "init_value"
$pump.Range("F1").Value = "hmi_tag"
$pump.Range("A1:F1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# Populate column-by-column (type, then init_value, then variable_name,
# then hmi_tag) so newly-introduced shared strings land in the same
# order the original generator produced them in.
# ---------------------------------------------------------------------

# Column D - type
$pump.Range("D2").Value = "ARRAY [3] OF BOOL"
$pump.Range("D3").Value = "WORD"
$pump.Range("D4").Value = "ARRAY [3] OF WORD"
$pump.Range("D5").Value = "ARRAY [3] OF WORD"
$pump.Range("D6").Value = "ARRAY [3] OF WORD"
$pump.Range("D7").Value = "ARRAY [4] OF BOOL"
$pump.Range("D8").Value = "WORD"
$pump.Range("D9").Value = "ARRAY [3] OF BOOL"
$pump.Range("D10").Value = "ARRAY [3] OF WORD"
$pump.Range("D11").Value = "ARRAY [3] OF WORD"
$pump.Range("D12").Value = "ARRAY [3] OF WORD"
$pump.Range("D13").Value = "ARRAY [4] OF BOOL"
$pump.Range("D2:D13").HorizontalAlignment = -4108

# Column E - init_value
$pump.Range("E2").Value = "[3(FALSE)]"
$pump.Range("E3").Value = 0
$pump.Range("E4").Value = "[3(0)]"
$pump.Range("E5").Value = "[3(0)]"
$pump.Range("E6").Value = "[3(0)]"
$pump.Range("E7").Value = "[4(FALSE)]"
$pump.Range("E8").Value = 0
$pump.Range("E9").Value = "[3(FALSE)]"
$pump.Range("E10").Value = "[3(0)]"
$pump.Range("E11").Value = "[3(0)]"
$pump.Range("E12").Value = "[3(0)]"
$pump.Range("E13").Value = "[4(FALSE)]"
$pump.Range("E9").HorizontalAlignment = -4108

# Column B - variable_name
$pump.Range("B2").Value = "pump_e_OnArr"
$pump.Range("B3").Value = "pump_e_FlowOrRpm"
$pump.Range("B4").Value = "pump_e_FlowSetArr"
$pump.Range("B5").Value = "pump_e_RpmSetArr"
$pump.Range("B6").Value = "pump_e_RpmValArr"
$pump.Range("B7").Value = "pump_e_ValveOnArr"
$pump.Range("B8").Value = "pump_h_Mode"
$pump.Range("B9").Value = "pump_h_OnArr"
$pump.Range("B10").Value = "pump_h_FlowSetArr"
$pump.Range("B11").Value = "pump_h_RpmSetArr"
$pump.Range("B12").Value = "pump_h_RpmValArr"
$pump.Range("B13").Value = "pump_h_ValveOnArr"

# Column F - hmi_tag (only rows 8-13, plus the styled-but-empty F2)
$pump.Range("F2").HorizontalAlignment = -4108
$pump.Range("F8").Value = "x"
$pump.Range("F9").Value = "x"
$pump.Range("F10").Value = "x"
$pump.Range("F11").Value = "x"
$pump.Range("F12").Value = "x"
$pump.Range("F13").Value = "x"

# Column A / C - base_addr / addr_offset (numeric, no shared strings)
$pump.Range("A2").Value = 11000
$pump.Range("C2").Value = 0
$pump.Range("C3").Value = 1
$pump.Range("C4").Value = 2
$pump.Range("C5").Value = 5
$pump.Range("C6").Value = 8
$pump.Range("C7").Value = 11
$pump.Range("C8").Value = 12
$pump.Range("C9").Value = 13
$pump.Range("C10").Value = 14
$pump.Range("C11").Value = 17
$pump.Range("C12").Value = 20
$pump.Range("C13").Value = 23

# ---------------------------------------------------------------------
# 3. Column widths (best-fit, matching the other generator sheets).
# ---------------------------------------------------------------------
$pump.Columns.Item("A").ColumnWidth = 10.140625
$pump.Columns.Item("B").ColumnWidth = 23.85546875
$pump.Columns.Item("C").ColumnWidth = 11.28515625
$pump.Columns.Item("D").ColumnWidth = 19
$pump.Columns.Item("E").ColumnWidth = 10

# ---------------------------------------------------------------------
# 4. View state: Pump becomes the selected/active tab, scrolled so
#    column B is the leftmost visible column, zoomed at 100%.
# ---------------------------------------------------------------------
$pump.Activate()
$pump.Range("D21").Select()
$excel.ActiveWindow.Zoom = 100
